# Apply the Feb 16, 2023 12:18 UTC "cryptos" symbol-list refresh:
# updated prices (D), 1-hour volume deltas (E), and the snapshot hour (G)
# from "11" to "12" for every data row (2-51). Columns are stored as text
# in the source sheet, so NumberFormat is forced to Text ("@") before each
# assignment to keep the exact original string formatting (trailing zeros,
# percent signs, sign, etc.) instead of having the value re-interpreted as
# a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "318.96" },
    @{ Cell = "E2"; Value = "5.82%" },
    @{ Cell = "G2"; Value = "12" },
    @{ Cell = "D3"; Value = "48.65" },
    @{ Cell = "E3"; Value = "14.43%" },
    @{ Cell = "G3"; Value = "12" },
    @{ Cell = "D4"; Value = "5.228" },
    @{ Cell = "E4"; Value = "4.34%" },
    @{ Cell = "G4"; Value = "12" },
    @{ Cell = "D5"; Value = "0.08059" },
    @{ Cell = "E5"; Value = "5.43%" },
    @{ Cell = "G5"; Value = "12" },
    @{ Cell = "D6"; Value = "4.587" },
    @{ Cell = "E6"; Value = "4.31%" },
    @{ Cell = "G6"; Value = "12" },
    @{ Cell = "D7"; Value = "1.643" },
    @{ Cell = "E7"; Value = "2.28%" },
    @{ Cell = "G7"; Value = "12" },
    @{ Cell = "D8"; Value = "1.212" },
    @{ Cell = "E8"; Value = "21.89%" },
    @{ Cell = "G8"; Value = "12" },
    @{ Cell = "E9"; Value = "6.89%" },
    @{ Cell = "G9"; Value = "12" },
    @{ Cell = "D10"; Value = "0.1930" },
    @{ Cell = "E10"; Value = "4.72%" },
    @{ Cell = "G10"; Value = "12" },
    @{ Cell = "D11"; Value = "0.09423" },
    @{ Cell = "E11"; Value = "4.21%" },
    @{ Cell = "G11"; Value = "12" },
    @{ Cell = "D12"; Value = "0.04604" },
    @{ Cell = "E12"; Value = "11.96%" },
    @{ Cell = "G12"; Value = "12" },
    @{ Cell = "E13"; Value = "0.37%" },
    @{ Cell = "G13"; Value = "12" },
    @{ Cell = "D14"; Value = "0.001328" },
    @{ Cell = "E14"; Value = "3.86%" },
    @{ Cell = "G14"; Value = "12" },
    @{ Cell = "D15"; Value = "0.04179" },
    @{ Cell = "E15"; Value = "1.12%" },
    @{ Cell = "G15"; Value = "12" },
    @{ Cell = "D16"; Value = "0.005841" },
    @{ Cell = "E16"; Value = "-0.38%" },
    @{ Cell = "G16"; Value = "12" },
    @{ Cell = "D17"; Value = "3.337" },
    @{ Cell = "E17"; Value = "-0.34%" },
    @{ Cell = "G17"; Value = "12" },
    @{ Cell = "D18"; Value = "2.428" },
    @{ Cell = "E18"; Value = "1.82%" },
    @{ Cell = "G18"; Value = "12" },
    @{ Cell = "D19"; Value = "0.3406" },
    @{ Cell = "E19"; Value = "2.06%" },
    @{ Cell = "G19"; Value = "12" },
    @{ Cell = "D20"; Value = "8.160" },
    @{ Cell = "E20"; Value = "-2.21%" },
    @{ Cell = "G20"; Value = "12" },
    @{ Cell = "D21"; Value = "0.1393" },
    @{ Cell = "E21"; Value = "3.59%" },
    @{ Cell = "G21"; Value = "12" },
    @{ Cell = "D22"; Value = "0.2910" },
    @{ Cell = "E22"; Value = "-11.38%" },
    @{ Cell = "G22"; Value = "12" },
    @{ Cell = "D23"; Value = "0.001312" },
    @{ Cell = "E23"; Value = "3.97%" },
    @{ Cell = "G23"; Value = "12" },
    @{ Cell = "D24"; Value = "0.004243" },
    @{ Cell = "E24"; Value = "7.21%" },
    @{ Cell = "G24"; Value = "12" },
    @{ Cell = "D25"; Value = "0.0001353" },
    @{ Cell = "E25"; Value = "1.01%" },
    @{ Cell = "G25"; Value = "12" },
    @{ Cell = "D26"; Value = "0.0003546" },
    @{ Cell = "E26"; Value = "-95.23%" },
    @{ Cell = "G26"; Value = "12" },
    @{ Cell = "G27"; Value = "12" },
    @{ Cell = "G28"; Value = "12" },
    @{ Cell = "G29"; Value = "12" },
    @{ Cell = "G30"; Value = "12" },
    @{ Cell = "G31"; Value = "12" },
    @{ Cell = "G32"; Value = "12" },
    @{ Cell = "G33"; Value = "12" },
    @{ Cell = "G34"; Value = "12" },
    @{ Cell = "G35"; Value = "12" },
    @{ Cell = "G36"; Value = "12" },
    @{ Cell = "G37"; Value = "12" },
    @{ Cell = "E38"; Value = "10.20%" },
    @{ Cell = "G38"; Value = "12" },
    @{ Cell = "D39"; Value = "0.05647" },
    @{ Cell = "E39"; Value = "7.04%" },
    @{ Cell = "G39"; Value = "12" },
    @{ Cell = "E40"; Value = "-2.24%" },
    @{ Cell = "G40"; Value = "12" },
    @{ Cell = "D41"; Value = "0.007955" },
    @{ Cell = "E41"; Value = "2.98%" },
    @{ Cell = "G41"; Value = "12" },
    @{ Cell = "D42"; Value = "0.1441" },
    @{ Cell = "E42"; Value = "7.24%" },
    @{ Cell = "G42"; Value = "12" },
    @{ Cell = "D43"; Value = "0.007702" },
    @{ Cell = "E43"; Value = "4.84%" },
    @{ Cell = "G43"; Value = "12" },
    @{ Cell = "D44"; Value = "0.008701" },
    @{ Cell = "E44"; Value = "19.53%" },
    @{ Cell = "G44"; Value = "12" },
    @{ Cell = "D45"; Value = "0.3505" },
    @{ Cell = "E45"; Value = "15.98%" },
    @{ Cell = "G45"; Value = "12" },
    @{ Cell = "D46"; Value = "0.00006917" },
    @{ Cell = "E46"; Value = "7.63%" },
    @{ Cell = "G46"; Value = "12" },
    @{ Cell = "D47"; Value = "0.00000000752" },
    @{ Cell = "E47"; Value = "0.99%" },
    @{ Cell = "G47"; Value = "12" },
    @{ Cell = "D48"; Value = "0.05485" },
    @{ Cell = "E48"; Value = "19.41%" },
    @{ Cell = "G48"; Value = "12" },
    @{ Cell = "D49"; Value = "0.004008" },
    @{ Cell = "E49"; Value = "-4.59%" },
    @{ Cell = "G49"; Value = "12" },
    @{ Cell = "D50"; Value = "0.00002105" },
    @{ Cell = "E50"; Value = "0.99%" },
    @{ Cell = "G50"; Value = "12" },
    @{ Cell = "D51"; Value = "0.0002005" },
    @{ Cell = "E51"; Value = "0.99%" },
    @{ Cell = "G51"; Value = "12" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}
